# Updates the cryptos price/volume list (Sheet1) to the latest scraped values.
# Numeric-looking text cells (prices such as "1.01", "19.60", etc.) must be forced
# to stay as text - otherwise Excel auto-converts them to real numbers and strips
# the original formatting (e.g. trailing zeros). We do that by temporarily
# switching the cell to a text number format, assigning the value, then resetting
# the cell style back to "Normal" so no stray style attribute is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.996.14'
$ws.Range("D3").Value = '1.641.36'
$ws.Range("E3").Value = '  -0.14%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.60%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.03'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.33%  '
$ws.Range("E7").Value = '  +0.54%  '
$ws.Range("E8").Value = '  -0.47%  '
$ws.Range("E9").Value = '  +0.21%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.60'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.77%  '
$ws.Range("E11").Value = '  +0.63%  '
$ws.Range("E12").Value = '  -0.06%  '
$ws.Range("D13").Value = '1.867.09'
$ws.Range("E13").Value = '  -0.21%  '
$ws.Range("D14").Value = '1.622.90'
$ws.Range("E14").Value = '  -0.53%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.545'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.04%  '
$ws.Range("D16").Value = '0.0₃0765'
$ws.Range("E16").Value = '  +0.08%  '
$ws.Range("E17").Value = '  -0.44%  '
$ws.Range("D18").Value = '25.896.17'
$ws.Range("E18").Value = '  -1.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.01'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.57%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '192.88'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.04%  '
$ws.Range("E21").Value = '  -1.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.94'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.11%  '
$ws.Range("E23").Value = '  -0.31%  '
$ws.Range("B24").Value = 'Stellar'
$ws.Range("C24").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.131'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.08%  '
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.58'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.42%  '
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.80'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.66%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.01'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.57%  '
$ws.Range("E28").Value = '  -0.16%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.55'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.26%  '
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0499'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.40%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.28'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.24'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("E34").Value = '  -3.46%  '
$ws.Range("E35").Value = '  +2.20%  '
$ws.Range("E36").Value = '  -1.00%  '
$ws.Range("D37").Value = '1.134.50'
$ws.Range("E37").Value = '  +0.23%  '
$ws.Range("E38").Value = '  -1.68%  '
$ws.Range("E39").Value = '  -1.02%  '
$ws.Range("E41").Value = '  +0.28%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '99.58'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.83%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.796'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.11%  '
$ws.Range("D44").Value = '1.776.14'
$ws.Range("E44").Value = '  -0.23%  '
$ws.Range("E45").Value = '  +1.45%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '56.68'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.24%  '
$ws.Range("E47").Value = '  +2.82%  '
$ws.Range("E48").Value = '  -1.28%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.72'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.03%  '
$ws.Range("E50").Value = '  -0.44%  '
$ws.Range("E51").Value = '  -0.64%  '
